$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.050.96'
$ws.Range("E2").Value = '  -2.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.161.00'
$ws.Range("E3").Value = '  -7.58%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.61'
$ws.Range("E5").Value = '  -2.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.95'
$ws.Range("E6").Value = '  -5.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.616'
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.158.15'
$ws.Range("E9").Value = '  -7.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  -5.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.54'
$ws.Range("E11").Value = '  -6.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.393'
$ws.Range("E12").Value = '  -4.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.709.50'
$ws.Range("E13").Value = '  -7.53%  '
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.07'
$ws.Range("E15").Value = '  -7.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.953.59'
$ws.Range("E16").Value = '  -2.45%  '
$ws.Range("E17").Value = '  -6.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.162.45'
$ws.Range("E18").Value = '  -7.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.73'
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.79'
$ws.Range("E20").Value = '  -7.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '356.90'
$ws.Range("E21").Value = '  -3.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.27'
$ws.Range("E22").Value = '  -4.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.17'
$ws.Range("E24").Value = '  -5.41%  '
$ws.Range("E25").Value = '  -6.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.306.24'
$ws.Range("E26").Value = '  -7.52%  '
$ws.Range("E27").Value = '  -7.96%  '
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("E29").Value = '  -1.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.96'
$ws.Range("E33").Value = '  -5.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.31'
$ws.Range("E34").Value = '  -8.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("E35").Value = '  -4.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.62'
$ws.Range("E36").Value = '  -6.19%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.51'
$ws.Range("E37").Value = '  -3.52%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.44'
$ws.Range("E38").Value = '  -6.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.835'
$ws.Range("E39").Value = '  -3.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.03'
$ws.Range("E41").Value = '  -5.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.671.11'
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("E43").Value = '  -6.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.18'
$ws.Range("E45").Value = '  -3.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.55'
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0655'
$ws.Range("E47").Value = '  -4.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.04'
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '322.67'
$ws.Range("E49").Value = '  -3.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0272'
$ws.Range("E50").Value = '  -4.72%  '
$ws.Range("E51").Value = '  -1.42%  '
